# quant Infrastructure Fund holdings sheet refresh ("updated data from quant engine").
#
# The engine now also classifies each holding's month-over-month trend, so a
# new "Status" column is inserted right after "Mutual Fund", and the monthly
# window rolls forward by one month (Jan/Dec/Oct instead of Jan/Dec/Nov) which
# also changes every QoQ figure (QoQ now compares Jan_2026 against Oct_2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at D ("Status"). Excel shifts the former D:H
# ("Jan_2026","Dec_2025","Nov_2025","MoM","QoQ") one column right to E:I,
# carrying the header formatting (style 1: bold/centered/bordered) with it -
# so Jan_2026/Dec_2025/MoM figures for every row need no further changes.
$ws.Columns("D:D").Insert()

# --- Header row: name the new column and roll the trailing month forward ---
$ws.Range("D1").Value() = "Status"
$ws.Range("G1").Value() = "Oct_2025"

# --- This month, "Adani Enterprises Limited Rights" fully rolled off and the
#     fund freshly entered Infosys Limited in its place (same row slot) ---
$ws.Range("A26").Value() = "INE009A01021"
$ws.Range("B26").Value() = "Infosys Limited"

# --- Per-holding trend label (Status) + refreshed Oct_2025 weight + QoQ ---
$holdingUpdates = @(
    @{ Row = 2; Status = "Adding Consistently"; Oct_2025 = 7.658724; QoQ = 2.049822999999999 },
    @{ Row = 3; Status = "Adding Consistently"; Oct_2025 = 9.365955; QoQ = 0.3411360000000005 },
    @{ Row = 4; Status = "Adding Consistently"; Oct_2025 = 8.919873; QoQ = 0.1326969999999985 },
    @{ Row = 5; Status = "Adding"; Oct_2025 = 6.36157; QoQ = -0.1794560000000001 },
    @{ Row = 6; Status = "Reducing Consistently"; Oct_2025 = 5.860319; QoQ = -0.7280129999999998 },
    @{ Row = 7; Status = "Adding"; Oct_2025 = 6.466866; QoQ = -1.499530999999999 },
    @{ Row = 8; Status = "Adding"; Oct_2025 = 5.598654; QoQ = -1.030678 },
    @{ Row = 9; Status = "Reducing"; Oct_2025 = 1.576079; QoQ = 2.931865 },
    @{ Row = 10; Status = "Reducing Consistently"; Oct_2025 = 3.968668; QoQ = -0.5910980000000001 },
    @{ Row = 11; Status = "Adding"; Oct_2025 = 3.264441; QoQ = -0.01706099999999999 },
    @{ Row = 12; Status = "Reducing Consistently"; Oct_2025 = 2.840387; QoQ = -0.0941080000000003 },
    @{ Row = 13; Status = "Reducing"; Oct_2025 = 1.051472; QoQ = 1.535209 },
    @{ Row = 14; Status = "Adding Consistently"; Oct_2025 = 2.045133; QoQ = 0.3992200000000001 },
    @{ Row = 15; Status = "Reducing"; Oct_2025 = 1.706299; QoQ = 0.1066370000000001 },
    @{ Row = 16; Status = "Adding Consistently"; Oct_2025 = 1.1804; QoQ = 0.117712 },
    @{ Row = 17; Status = "Reducing Consistently"; Oct_2025 = 1.319649; QoQ = -0.140487 },
    @{ Row = 18; Status = "Reducing Consistently"; Oct_2025 = 1.118824; QoQ = -0.008216000000000001 },
    @{ Row = 19; Status = "Reducing Consistently"; Oct_2025 = 1.308332; QoQ = -0.2512990000000002 },
    @{ Row = 20; Status = "Adding Consistently"; Oct_2025 = 0.58827; QoQ = 0.459062 },
    @{ Row = 21; Status = "Adding"; Oct_2025 = 8.313412; QoQ = -7.385738999999999 },
    @{ Row = 22; Status = "Reducing Consistently"; Oct_2025 = 0.783232; QoQ = -0.002846000000000015 },
    @{ Row = 23; Status = "Fresh Entry"; Oct_2025 = 0; QoQ = 0.686119 },
    @{ Row = 24; Status = "Adding"; Oct_2025 = 0.488493; QoQ = -0.3414970000000001 },
    @{ Row = 25; Status = "Complete Exit"; Oct_2025 = 2.776239; QoQ = -2.776239 },
    @{ Row = 26; Status = "Complete Exit"; Oct_2025 = 2.768339; QoQ = -2.768339 },
    @{ Row = 27; Status = "Complete Exit"; Oct_2025 = 7.524506; QoQ = -7.524506 }
)

foreach ($update in $holdingUpdates) {
    $ws.Cells.Item($update.Row, 4).Value() = $update.Status
    $ws.Cells.Item($update.Row, 7).Value() = $update.Oct_2025
    $ws.Cells.Item($update.Row, 9).Value() = $update.QoQ
}
